# Purchases_Totals.xlsx — fill in Robert's purchases with his itemized
# receipts, which ripples through the Totals sheet formulas.

$wb = $excel.ActiveWorkbook

# --- "Robert's purchases" sheet ---
$ws = $wb.Worksheets.Item("Robert's purchases")
$ws.Activate()

# The previous "even split" placeholder had E5:E7 merged together; once
# real per-row amounts go in, that merge (and its special centered style)
# needs to go so every row can carry its own value.
$ws.Range("E5:E7").UnMerge()
$ws.Range("E5:E7").ClearFormats()
$ws.Range("E5:E7").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# Component names — entered row 5 up to row 2, then row 6 down to row 9,
# matching how the shared-string table ends up ordered.
$ws.Range("A5").Value = "Capacitors"
$ws.Range("A4").Value = "OSHpark"
$ws.Range("A3").Value = "Digikey 2"
$ws.Range("A2").Value = "Digikey 1"
$ws.Range("A6").Value = "Magnets"
$ws.Range("A7").Value = "IR Sensor"
$ws.Range("A8").Value = "Battery Charger"
$ws.Range("A9").Value = "Battery"

# Per-item prices.
$ws.Range("E3").Value = 23.82
$ws.Range("E4").Value = 21.25
$ws.Range("E5").Value = 15.94
$ws.Range("E6").Value = 6.41
$ws.Range("E7").Value = 12.71
$ws.Range("E8").Value = 7.48
$ws.Range("E9").Value = 20.52

# First item (row 2) got pasted in from a different source with its own
# font formatting (Arial, dark grey) instead of the sheet's usual font.
$ws.Range("E2").Font.Name = "Arial"
$ws.Range("E2").Font.Color = 2236962
$ws.Range("E2").Value = 74.53

# Total now sums the whole item list instead of sitting empty.
$ws.Range("E10").Formula = "=SUM(E2:E9)"

# A couple of stray formatted-but-empty cells far to the right (K8, Q8) —
# left over from when this was pasted in from a wider sheet.
$ws.Range("K8").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("Q8").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

$ws.Range("E11").Select()

# --- "Totals" sheet: recompute once Robert's column has real numbers ---
$totals = $wb.Worksheets.Item("Totals")
$totals.Activate()
$totals.Range("H30").Select()
